$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date column (A2:A6) to the new menu date
$ws.Range("A2:A6").Value = "15.05.2025"

# Update the dish descriptions in column B
$ws.Range("B2").Value = "Panierowana karkówka, ziemniaki, surówka + zupa ogórkowa lub kapuśniak ze świeżej kapusty"
$ws.Range("B3").Value = "Kawałki kurczaka w sosie pieczarkowym, ziemniaki, surówka + zupa ogórkowa lub kapuśniak ze świeżej kapusty"
$ws.Range("B4").Value = "Panierowana karkówka, ziemniaki, surówka "
$ws.Range("B5").Value = "Kawałki kurczaka w sosie pieczarkowym, ziemniaki, surówka"
$ws.Range("B6").Value = "Zupa ogórkowa lub kapuśniak ze świeżej kapusty"

# Update prices in column C
$ws.Range("C2").Value = 31
$ws.Range("C3").Value = 31
$ws.Range("C4").Value = 28
$ws.Range("C5").Value = 28

# Update the active cell selection
$ws.Range("B8").Select()
